# Re-pulled data: update dSF (column F) values for several rows to match
# the freshly pulled source data / recalculated means.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -8
    8  = 2
    9  = -5
    18 = -7
    19 = -1
    22 = 0
    23 = 0
    28 = -1
    29 = -2
    32 = -2
    34 = -5
    35 = 0
    38 = -7
    40 = -4
    41 = 0
    42 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
